# Aula 16/02 — integração com redes, deploy
#
# The paragraph that ends in "Injeção de dependência " (which also carries
# the trailing _GoBack bookmark) is split into three paragraphs:
#   1. The original run "Injeção de dependência " (paragraph now ends here)
#   2. A brand-new, completely empty paragraph
#   3. A new paragraph that keeps the _GoBack bookmark and adds the
#      "DEPLOY:" heading plus the "Quem criar um servidor de API já cria
#      o de banco" note beneath it.

$d = $word.ActiveDocument

# --- Locate the target paragraph -------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute("Injeção de dependência ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Injeção de dependência' paragraph"
}
$para = $findRng.Paragraphs.First
$target = $para.Range

# Recover the real opening <w:p ...> tag (rsid attributes etc.) for the first
# of the three resulting paragraphs so it is left byte-identical to the
# original - strip any w14:paraId/w14:textId the round-trip getter may add,
# those are not present in the stored part.
$openXml = $target.WordOpenXML
$pOpenTag = '<w:p>'
if ($openXml -match '(<w:p [^>]*>)') {
    $pOpenTag = $matches[1]
    $pOpenTag = $pOpenTag -replace ' w14:paraId="[^"]*"', ''
    $pOpenTag = $pOpenTag -replace ' w14:textId="[^"]*"', ''
}

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Shared formatting fragments --------------------------------------------------
$pPrCommon = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="160" w:afterAutospacing="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="000000"/></w:rPr></w:pPr>'
$rPrBoldBlack = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="000000"/></w:rPr>'
$rPrPasted = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:szCs w:val="22"/><w:lang w:eastAsia="en-US"/></w:rPr>'

# --- Paragraph 1: unchanged run, now its own paragraph ----------------------------
$p1 = $pOpenTag + $pPrCommon + '<w:r>' + $rPrBoldBlack + '<w:t xml:space="preserve">Injeção de dependência </w:t></w:r></w:p>'

# --- Paragraph 2: brand-new, entirely empty ----------------------------------------
$p2 = '<w:p ' + $wns + '>' + $pPrCommon + '</w:p>'

# --- Paragraph 3: keeps the bookmark, gets the new DEPLOY note --------------------
$runs  = '<w:r>' + $rPrBoldBlack + '<w:t>DEPLOY:</w:t></w:r>'
$runs += '<w:r>' + $rPrBoldBlack + '<w:br/></w:r>'
$runs += '<w:r>' + $rPrPasted + '<w:t>Q</w:t></w:r>'
$runs += '<w:r>' + $rPrPasted + '<w:t xml:space="preserve">uem criar um servidor de </w:t></w:r>'
$runs += '<w:r>' + $rPrPasted + '<w:t>API</w:t></w:r>'
$runs += '<w:r>' + $rPrPasted + '<w:t xml:space="preserve"> </w:t></w:r>'
$runs += '<w:r>' + $rPrPasted + '<w:t>já</w:t></w:r>'
$runs += '<w:r>' + $rPrPasted + '<w:t xml:space="preserve"> cria o de banco</w:t></w:r>'

$p3 = '<w:p ' + $wns + '>' + $pPrCommon + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + $runs + '</w:p>'

$fragment = $p1 + $p2 + $p3

$target.InsertXML($fragment)
